$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 21: MI vs SRH (contest 9) - points values
$ws.Range("E21").Value = 30
$ws.Range("H21").Value = 80
$ws.Range("K21").Value = 40
$ws.Range("N21").Value = 20
$ws.Range("Q21").Value = 50
$ws.Range("T21").Value = 0
$ws.Range("W21").Value = 100
$ws.Range("Z21").Value = 70
$ws.Range("AC21").Value = 60

# Row 22: RCB vs KKR (contest 10) - points values
$ws.Range("E22").Value = 100
$ws.Range("H22").Value = 0
$ws.Range("K22").Value = 70
$ws.Range("N22").Value = 30
$ws.Range("Q22").Value = 50
$ws.Range("T22").Value = 20
$ws.Range("W22").Value = 60
$ws.Range("Z22").Value = 40
$ws.Range("AC22").Value = 80

# Row 23: DC vs PBKS (contest 11) - points values
$ws.Range("E23").Value = 80
$ws.Range("H23").Value = 0
$ws.Range("K23").Value = 20
$ws.Range("N23").Value = 70
$ws.Range("Q23").Value = 30
$ws.Range("T23").Value = 60
$ws.Range("W23").Value = 40
$ws.Range("Z23").Value = 100
$ws.Range("AC23").Value = 50
